# The edit moves the "_GoBack" bookmark (an empty bookmark Word uses to
# remember the last edit location) from the end of the paragraph that ends
# in "//" to the end of the following paragraph, "后续优化：座位画成椅子"
# (right after its text, before the paragraph mark).
#
# Word keeps "_GoBack" unique, so adding it again at the new location
# automatically removes it from its old location.

$d = $word.ActiveDocument

$found = $d.Content
$found.Find.Execute("后续优化：座位画成椅子", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0)

if ($found.Find.Found) {
    # Collapse the found range to its end -- right after the last
    # character of the paragraph's text, immediately before the paragraph
    # mark. That is exactly where "_GoBack" needs to end up.
    $target = $d.Range($found.End, $found.End)

    # A zero-length Range that already sits exactly on a paragraph-mark
    # boundary cannot be handed straight to Bookmarks.Add -- doing so
    # resets the bookmark back to the start of the document. Working
    # around it by briefly inserting a one-character placeholder,
    # bookmarking that character, and then deleting it through the
    # bookmark's own Range collapses the bookmark back to zero length
    # exactly where it needs to be, without hitting that problem.
    $target.InsertAfter("X")
    $d.Bookmarks.Add("_GoBack", $target)
    $bm = $d.Bookmarks("_GoBack")
    $bm.Range.Text = ""
}
